$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per latest refresh
$ws.Range('D2').Value = '42.497.75'
$ws.Range('E2').Value = '  -0.99%  '
$ws.Range('D3').Value = '2.254.37'
$ws.Range('E3').Value = '  -1.39%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '''248.88'
$ws.Range('E5').Value = '  -1.31%  '
$ws.Range('D6').Value = '''0.634'
$ws.Range('E6').Value = '  +0.90%  '
$ws.Range('D7').Value = '''76.91'
$ws.Range('E7').Value = '  +6.58%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = '''0.642'
$ws.Range('E9').Value = '  -0.12%  '
$ws.Range('D10').Value = '''40.47'
$ws.Range('E10').Value = '  +3.87%  '
$ws.Range('D11').Value = '''0.0958'
$ws.Range('E11').Value = '  -2.07%  '
$ws.Range('D12').Value = '''7.30'
$ws.Range('E12').Value = '  -1.28%  '
$ws.Range('E13').Value = '  -0.26%  '
$ws.Range('D14').Value = '2.595.80'
$ws.Range('E14').Value = '  -1.32%  '
$ws.Range('D15').Value = '''14.96'
$ws.Range('E15').Value = '  -1.08%  '
$ws.Range('D16').Value = '''0.864'
$ws.Range('E16').Value = '  -2.51%  '
$ws.Range('D17').Value = '2.265.91'
$ws.Range('E17').Value = '  -1.24%  '
$ws.Range('D18').Value = '42.428.17'
$ws.Range('E18').Value = '  -1.15%  '
$ws.Range('D19').Value = '0.0₃0986'
$ws.Range('E19').Value = '  -1.87%  '
$ws.Range('E20').Value = '  -2.45%  '
$ws.Range('E21').Value = '  -2.35%  '
$ws.Range('D22').Value = '''232.33'
$ws.Range('E22').Value = '  -1.69%  '
$ws.Range('D23').Value = '''2.16'
$ws.Range('E23').Value = '  -4.90%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').Value = '''3.74'
$ws.Range('E25').Value = '  -5.34%  '
$ws.Range('D26').Value = '''11.26'
$ws.Range('E26').Value = '  -3.15%  '
$ws.Range('E27').Value = '  -5.56%  '
$ws.Range('E28').Value = '  -0.71%  '
$ws.Range('D29').Value = '''169.26'
$ws.Range('E29').Value = '  +0.97%  '
$ws.Range('D30').Value = '''6.83'
$ws.Range('E30').Value = '  +8.24%  '
$ws.Range('D31').Value = '''20.64'
$ws.Range('E31').Value = '  -2.38%  '
$ws.Range('D32').Value = '''0.0853'
$ws.Range('E32').Value = '  +6.24%  '
$ws.Range('D34').Value = '''30.75'
$ws.Range('E34').Value = '  -0.99%  '
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('D36').Value = '''4.50'
$ws.Range('E36').Value = '  -4.98%  '
$ws.Range('D37').Value = '''4.74'
$ws.Range('E37').Value = '  -0.97%  '
$ws.Range('D38').Value = '''0.0300'
$ws.Range('E38').Value = '  -3.17%  '
$ws.Range('D39').Value = '''13.13'
$ws.Range('E39').Value = '  -1.91%  '
$ws.Range('E40').Value = '  -4.14%  '
$ws.Range('E41').Value = '  -1.79%  '
$ws.Range('D42').Value = '''118.80'
$ws.Range('E42').Value = '  +23.13%  '
$ws.Range('E43').Value = '  -2.20%  '
$ws.Range('E44').Value = '  -2.65%  '
$ws.Range('D45').Value = '''8.85'
$ws.Range('E45').Value = '  -4.54%  '
$ws.Range('E46').Value = '  -2.41%  '
$ws.Range('E47').Value = '  -0.46%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').Value = '''1.14'
$ws.Range('E48').Value = '  -3.85%  '
$ws.Range('B49').Value = 'FTXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D49').Value = '''4.41'
$ws.Range('E49').Value = '  -10.68%  '
$ws.Range('D50').Value = '''1.18'
$ws.Range('E50').Value = '  -1.33%  '
$ws.Range('D51').Value = '''4.14'
$ws.Range('E51').Value = '  -2.53%  '
